$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 117 and 118: swap all data (columns B..AC), keep the sequential
# "A" index (115 / 116) fixed in place ---

# New row 117 (was row 118's data)
$ws.Range("B117").Value = 7013702
$ws.Range("F117").Value = "Defensor Sporting"
$ws.Range("G117").Value = "Danubio"
$ws.Range("I117").Value = 2
$ws.Range("K117").Value = 1.8
$ws.Range("L117").Value = 3.6
$ws.Range("M117").Value = 4.2
$ws.Range("N117").Value = 1.8
$ws.Range("O117").Value = 3.6
$ws.Range("P117").Value = 4.2
$ws.Range("Q117").Value = -0.75
$ws.Range("R117").Value = 2.05
$ws.Range("S117").Value = 1.8
$ws.Range("T117").Value = 2.25
$ws.Range("U117").Value = 1.85
$ws.Range("V117").Value = 2
$ws.Range("Y117").Value = 3.2
$ws.Range("AA117").Value = 0.8
$ws.Range("AB117").Value = -0.5
$ws.Range("AC117").Value = 0.5

# New row 118 (was row 117's data)
$ws.Range("B118").Value = 7013886
$ws.Range("F118").Value = "Racing Club de Montevideo"
$ws.Range("G118").Value = "Cerro"
$ws.Range("I118").Value = 1
$ws.Range("K118").Value = 2.25
$ws.Range("L118").Value = 3.1
$ws.Range("M118").Value = 3.25
$ws.Range("N118").Value = 2.25
$ws.Range("O118").Value = 2.875
$ws.Range("P118").Value = 3.5
$ws.Range("Q118").Value = -0.25
$ws.Range("R118").Value = 1.95
$ws.Range("S118").Value = 1.9
$ws.Range("T118").Value = 2
$ws.Range("U118").Value = 1.925
$ws.Range("V118").Value = 1.925
$ws.Range("Y118").Value = 2.5
$ws.Range("AA118").Value = 0.8999999999999999
$ws.Range("AB118").Value = -1
$ws.Range("AC118").Value = 0.925

# --- Row 187: updated odds from a later snapshot ---
$ws.Range("B187").Value = 8081251
$ws.Range("E187").Value = 45403.54166666666
$ws.Range("F187").Value = "Atletico Fenix Montevideo"
$ws.Range("G187").Value = "Montevideo Wanderers"
$ws.Range("K187").Value = 2.5
$ws.Range("L187").Value = 3
$ws.Range("M187").Value = 3
$ws.Range("N187").Value = 2.5
$ws.Range("O187").Value = 2.9
$ws.Range("P187").Value = 3.1
$ws.Range("R187").Value = 2.1
$ws.Range("S187").Value = 1.775
$ws.Range("T187").Value = 2
$ws.Range("U187").Value = 1.875
$ws.Range("V187").Value = 1.975

# --- Row 188: updated odds from a later snapshot ---
$ws.Range("B188").Value = 8081885
$ws.Range("E188").Value = 45403.64583333334
$ws.Range("F188").Value = "CA River Plate"
$ws.Range("G188").Value = "Club Atletico Progreso"
$ws.Range("K188").Value = 2.625
$ws.Range("L188").Value = 3.1
$ws.Range("M188").Value = 2.75
$ws.Range("N188").Value = 3.1
$ws.Range("O188").Value = 3.1
$ws.Range("P188").Value = 2.375
$ws.Range("Q188").Value = 0.25
$ws.Range("R188").Value = 1.8
$ws.Range("S188").Value = 2.05
$ws.Range("U188").Value = 1.975
$ws.Range("V188").Value = 1.875

# --- Rows 189 and 190 were removed entirely from the source data ---
$ws.Rows("189:190").Delete()
